$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-08 05:25:28"
$wsDeDe.Range("H2").Value = "2016-09-08 05:25:28"

$wsZhCn.Range("H2").Value = "2016-09-08 05:25:23"
$wsZhCn.Range("K2").Value = "2016-09-08 05:25:48"

$wsDeDe.Range("K2").Value = "2016-09-08 05:25:57"
